# Auto-generated script to apply 2024-03-21 violent crime data update
$wb = $excel.ActiveWorkbook

# --- Sheet: Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1539  # Aggravated Assault: 1515 -> 1539
$ws.Range("K3").Value = 1467  # Aggravated Battery: 1452 -> 1467
$ws.Range("K5").Value = 97  # Homicide: 96 -> 97
$ws.Range("K6").Value = 1884  # Robbery: 1862 -> 1884
$ws.Range("K7").Value = 5298  # Total: 5236 -> 5298

# --- Sheet: Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 100  # Aggravated Assault: 99 -> 100
$ws.Range("K3").Value = 98  # Aggravated Battery: 97 -> 98
$ws.Range("K6").Value = 110  # Robbery: 109 -> 110
$ws.Range("K7").Value = 332  # Total: 329 -> 332

# --- Sheet: South Chicago ---
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 47  # Aggravated Assault: 46 -> 47
$ws.Range("K3").Value = 36  # Aggravated Battery: 34 -> 36
$ws.Range("K7").Value = 107  # Total: 104 -> 107

# --- Sheet: Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 63  # Aggravated Assault: 59 -> 63
$ws.Range("K4").Value = 14  # Criminal Sexual Assault: 13 -> 14
$ws.Range("K6").Value = 55  # Robbery: 53 -> 55
$ws.Range("K7").Value = 216  # Total: 209 -> 216

# --- Sheet: West Pullman ---
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 26  # Aggravated Battery: 25 -> 26
$ws.Range("K7").Value = 90  # Total: 89 -> 90

# --- Sheet: Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 59  # Aggravated Battery: 58 -> 59
$ws.Range("K7").Value = 178  # Total: 177 -> 178

# --- Sheet: New City ---
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 37  # Aggravated Assault: 35 -> 37
$ws.Range("K4").Value = 3  # Criminal Sexual Assault: 2 -> 3
$ws.Range("K6").Value = 57  # Robbery: 56 -> 57
$ws.Range("K7").Value = 133  # Total: 129 -> 133

# --- Sheet: By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 38  # Albany Park: 37 -> 38
$ws.Range("K4").Value = 22  # Archer Heights: 21 -> 22
$ws.Range("K7").Value = 151  # Auburn Gresham: 148 -> 151
$ws.Range("K8").Value = 332  # Austin: 329 -> 332
$ws.Range("K9").Value = 27  # Avalon Park: 26 -> 27
$ws.Range("K18").Value = 42  # Calumet Heights: 40 -> 42
$ws.Range("K19").Value = 135  # Chatham: 134 -> 135
$ws.Range("K21").Value = 17  # Chinatown: 16 -> 17
$ws.Range("K23").Value = 55  # Douglas: 53 -> 55
$ws.Range("K29").Value = 241  # Englewood: 236 -> 241
$ws.Range("K31").Value = 63  # Gage Park: 61 -> 63
$ws.Range("K33").Value = 216  # Garfield Park: 209 -> 216
$ws.Range("K36").Value = 58  # Grand Boulevard: 57 -> 58
$ws.Range("K37").Value = 178  # Grand Crossing: 177 -> 178
$ws.Range("K41").Value = 56  # Hermosa: 55 -> 56
$ws.Range("K42").Value = 182  # Humboldt Park: 181 -> 182
$ws.Range("K43").Value = 52  # Hyde Park: 51 -> 52
$ws.Range("K44").Value = 49  # Irving Park: 48 -> 49
$ws.Range("K48").Value = 57  # Lake View: 56 -> 57
$ws.Range("K51").Value = 64  # Little Italy, UIC: 61 -> 64
$ws.Range("K55").Value = 56  # Lower West Side: 55 -> 56
$ws.Range("K60").Value = 42  # Morgan Park: 41 -> 42
$ws.Range("K63").Value = 17  # NO NEIGHBORHOOD DATA: 18 -> 17
$ws.Range("K65").Value = 133  # New City: 129 -> 133
$ws.Range("K67").Value = 205  # North Lawndale: 203 -> 205
$ws.Range("K76").Value = 72  # River North: 71 -> 72
$ws.Range("K78").Value = 76  # Rogers Park: 74 -> 76
$ws.Range("K79").Value = 143  # Roseland: 140 -> 143
$ws.Range("K83").Value = 107  # South Chicago: 104 -> 107
$ws.Range("K88").Value = 68  # United Center: 66 -> 68
$ws.Range("K90").Value = 51  # Washington Heights: 50 -> 51
$ws.Range("K91").Value = 52  # Washington Park: 51 -> 52
$ws.Range("K94").Value = 66  # West Loop: 65 -> 66
$ws.Range("K95").Value = 90  # West Pullman: 89 -> 90
$ws.Range("K96").Value = 72  # West Ridge: 71 -> 72
$ws.Range("K101").Value = 5298  # Total: 5236 -> 5298

# --- Sheet: Gage Park ---
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 25  # Aggravated Assault: 24 -> 25
$ws.Range("K3").Value = 10  # Aggravated Battery: 9 -> 10
$ws.Range("K7").Value = 63  # Total: 61 -> 63

# --- Sheet: North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 60  # Aggravated Assault: 58 -> 60
$ws.Range("K7").Value = 205  # Total: 203 -> 205

# --- Sheet: Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 63  # Aggravated Assault: 62 -> 63
$ws.Range("K3").Value = 80  # Aggravated Battery: 79 -> 80
$ws.Range("K6").Value = 81  # Robbery: 78 -> 81
$ws.Range("K7").Value = 241  # Total: 236 -> 241

# --- Sheet: Lake View ---
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 25  # Robbery: 24 -> 25
$ws.Range("K7").Value = 57  # Total: 56 -> 57

# --- Sheet: Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K5").Value = 8  # Homicide: 7 -> 8
$ws.Range("K7").Value = 135  # Total: 134 -> 135

# --- Sheet: Irving Park ---
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 20  # Robbery: 19 -> 20
$ws.Range("K7").Value = 49  # Total: 48 -> 49

# --- Sheet: River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 40  # Robbery: 39 -> 40
$ws.Range("K7").Value = 72  # Total: 71 -> 72

# --- Sheet: Hermosa ---
$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 17  # Aggravated Assault: 16 -> 17
$ws.Range("K7").Value = 56  # Total: 55 -> 56

# --- Sheet: Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 43  # Aggravated Assault: 42 -> 43
$ws.Range("K4").Value = 9  # Criminal Sexual Assault: 10 -> 9
$ws.Range("K6").Value = 80  # Robbery: 79 -> 80
$ws.Range("K7").Value = 182  # Total: 181 -> 182

# --- Sheet: Rogers Park ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 24  # Aggravated Assault: 23 -> 24
$ws.Range("K6").Value = 25  # Robbery: 24 -> 25
$ws.Range("K7").Value = 76  # Total: 74 -> 76

# --- Sheet: Lower West Side ---
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 22  # Aggravated Assault: 21 -> 22
$ws.Range("K7").Value = 56  # Total: 55 -> 56

# --- Sheet: Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 13  # Robbery: 11 -> 13
$ws.Range("K7").Value = 55  # Total: 53 -> 55

# --- Sheet: West Ridge ---
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 36  # Robbery: 35 -> 36
$ws.Range("K7").Value = 72  # Total: 71 -> 72

# --- Sheet: Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K6").Value = 11  # Robbery: 10 -> 11
$ws.Range("K7").Value = 52  # Total: 51 -> 52

# --- Sheet: Chinatown ---
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K6").Value = 10  # Robbery: 9 -> 10
$ws.Range("K7").Value = 17  # Total: 16 -> 17

# --- Sheet: Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 52  # Aggravated Assault: 51 -> 52
$ws.Range("K3").Value = 48  # Aggravated Battery: 47 -> 48
$ws.Range("K6").Value = 30  # Robbery: 29 -> 30
$ws.Range("K7").Value = 143  # Total: 140 -> 143

# --- Sheet: Calumet Heights ---
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 11  # Aggravated Battery: 9 -> 11
$ws.Range("K7").Value = 42  # Total: 40 -> 42

# --- Sheet: Grand Boulevard ---
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 23  # Aggravated Assault: 22 -> 23
$ws.Range("K7").Value = 58  # Total: 57 -> 58

# --- Sheet: Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 49  # Aggravated Battery: 47 -> 49
$ws.Range("K6").Value = 36  # Robbery: 35 -> 36
$ws.Range("K7").Value = 151  # Total: 148 -> 151

# --- Sheet: West Loop ---
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 29  # Robbery: 28 -> 29
$ws.Range("K7").Value = 66  # Total: 65 -> 66

# --- Sheet: Avalon Park ---
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 11  # Aggravated Battery: 10 -> 11
$ws.Range("K7").Value = 27  # Total: 26 -> 27

# --- Sheet: Albany Park ---
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 10  # Aggravated Assault: 9 -> 10
$ws.Range("K7").Value = 38  # Total: 37 -> 38

# --- Sheet: United Center ---
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 16  # Aggravated Assault: 15 -> 16
$ws.Range("K6").Value = 38  # Robbery: 37 -> 38
$ws.Range("K7").Value = 68  # Total: 66 -> 68

# --- Sheet: Washington Heights ---
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 12  # Robbery: 11 -> 12
$ws.Range("K7").Value = 51  # Total: 50 -> 51

# --- Sheet: Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 16  # Aggravated Assault: 14 -> 16
$ws.Range("K3").Value = 18  # Aggravated Battery: 17 -> 18
$ws.Range("K7").Value = 64  # Total: 61 -> 64

# --- Sheet: Morgan Park ---
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 13  # Robbery: 12 -> 13
$ws.Range("K7").Value = 42  # Total: 41 -> 42

# --- Sheet: Hyde Park ---
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 9  # Aggravated Assault: 8 -> 9
$ws.Range("K7").Value = 52  # Total: 51 -> 52

# --- Sheet: Archer Heights ---
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K2").Value = 7  # Aggravated Assault: 6 -> 7
$ws.Range("K6").Value = 22  # Total: 21 -> 22
